$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168, shifting existing rows 168..195 down to 169..196.
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(168, 1).Value = 5
$ws.Cells.Item(168, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(168, 3).Value = "Maule"
$ws.Cells.Item(168, 4).Value = 44505
$ws.Cells.Item(168, 5).Value = 7
$ws.Cells.Item(168, 6).Value = 100112009
$ws.Cells.Item(168, 7).Value = "Acelga"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 500
$ws.Cells.Item(168, 11).Value = 2000
$ws.Cells.Item(168, 12).Value = 2000
$ws.Cells.Item(168, 13).Value = 2000
$ws.Cells.Item(168, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 500
$ws.Cells.Item(168, 17).Value = 4
$ws.Cells.Item(168, 18).Value = "Hortaliza"
